# #456 Add test case
# Adds a new worksheet "header rows test" after the existing "Object测试"
# sheet. The new sheet contains a header row starting at row 3 (instead of
# row 1) followed by two sample data rows, mirroring the first two data
# rows of the original sheet - used to exercise "header rows" support.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "header rows test"

# ---- Header row (row 3) ----
$ws2.Range("A3").Value = "渠道ID"
$ws2.Range("B3").Value = "游戏"
$ws2.Range("C3").Value = "账号"
$ws2.Range("D3").Value = "注册时间"
$ws2.Range("E3").Value = "是否满30级"
$ws2.Range("F3").Value = "VIP"

$ws2.Range("C3").Font.Name = "宋体"
$ws2.Range("C3").VerticalAlignment = -4108

# ---- Data row 4 ----
$ws2.Range("A4").Value = 4
$ws2.Range("B4").Value = "极品飞车"
$ws2.Range("C4").Value = "XuSu2gFg32"
$ws2.Range("D4").Value = 43425
$ws2.Range("E4").Value = $true
$ws2.Range("F4").Value = "F"

# ---- Data row 5 ----
$ws2.Range("A5").Value = 8
$ws2.Range("B5").Value = "怪物世界"
$ws2.Range("C5").Value = "kxwWgaB"
$ws2.Range("D5").Value = 43425
$ws2.Range("E5").Value = $true
$ws2.Range("F5").Value = "N"

$dataRange = $ws2.Range("A4:F5")
$dataRange.Font.Name = "宋体"
$dataRange.VerticalAlignment = -4108

$numCol = $ws2.Range("A4:A5")
$numCol.NumberFormat = "0"
$numCol.HorizontalAlignment = -4152

$textCols = $ws2.Range("B4:C5")
$textCols.HorizontalAlignment = -4131

$dateCol = $ws2.Range("D4:D5")
$dateCol.NumberFormat = "yyyy-mm-dd"
$dateCol.HorizontalAlignment = -4108

$flagCols = $ws2.Range("E4:F5")
$flagCols.HorizontalAlignment = -4108


$ws2.Range("A3").Select() | Out-Null
$ws1.Activate() | Out-Null
